$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 37
$ws.Range("D2").Value = 3
$ws.Range("B5").Value = 0.925
$ws.Range("D5").Value = 0.075
